$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in new row 7 data (columns B-K) to accompany the existing A7 "testb1" label,
# mirroring the other test rows (collision rate, avg/max probe depth, std deviation,
# mode, median, collisions, sparsity, dims, nbuckets)
$ws.Range("B7").Value = 0.99519999999999997
$ws.Range("B7").NumberFormat = "0.00%"

$ws.Range("C7").Value = 210.5
$ws.Range("D7").Value = 470
$ws.Range("E7").Value = 77.7
$ws.Range("F7").Value = 156
$ws.Range("G7").Value = 1556.5
$ws.Range("H7").Value = 39810

$ws.Range("I7").Value = 0.99
$ws.Range("I7").NumberFormat = "0.00%"

$ws.Range("J7").Value = 20000
$ws.Range("K7").Value = 65536

# Update the selected/active cell to E12
$ws.Range("E12").Select()
